$d = $word.ActiveDocument

# Hunk 1 & 3: "{/}{/}{#Extraccions.length!=0}" -> "{/}{/}{/}{#Extraccions.length!=0}"
# (adds one extra "{/}" closing tag before the Extraccions.length!=0 block)
$d.Content.Find.Execute(
    "{/}{/}{#Extraccions.length!=0}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{/}{/}{/}{#Extraccions.length!=0}", 2
)

# Hunk 2 & 4: "{/Extraccions}{/}{/}{/}" -> "{/Extraccions}{/}{/}"
# (removes one trailing "{/}" closing tag)
$d.Content.Find.Execute(
    "{/Extraccions}{/}{/}{/}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{/Extraccions}{/}{/}", 2
)
